# 自动更新Excel文件 - 2025-10-19 23:11:28
#
# Refreshes the "剩余" (days remaining) / "开始时间" (start date) tracking
# columns (E / F) against "today" = 2025-10-20:
#   - If a row's cycle has run out (today has reached/passed start+total
#     days), the cycle restarts: 开始时间(F) is set to today and
#     剩余(E) is reset to the full 总天(D) count.
#   - Otherwise the row is still mid-cycle, so 开始时间(F) stays put and
#     剩余(E) simply ticks down by one day.
# Row 36 is intentionally left alone - its 开始时间 value is a malformed
# date ("202510929") so it can't be evaluated and is skipped, same as the
# source update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose cycle expired: restart (E = total days, F = today)
$resetRows = @(2,3,4,6,8,10,11,13,14,15,16,40,41,44,46,48,70,71,72,73,74,75,76,87,88,89,90,92,96,97,98,99)
$today = 20251020

foreach ($r in $resetRows) {
  $totalDays = $ws.Range("D" + $r).Value2
  $ws.Range("E" + $r).Value = $totalDays
  $ws.Range("F" + $r).Value = $today
}

# Rows still within their cycle: just decrement remaining days by one
$decrementRows = @(5,7,9,12,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,42,43,45,47,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,77,78,79,80,81,82,83,84,85,86,91,93,94,95)

foreach ($r in $decrementRows) {
  $remaining = $ws.Range("E" + $r).Value2
  $ws.Range("E" + $r).Value = $remaining - 1
}
